# Fixed a bug in holdAndWin: correct the row data in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (rows 2-21, columns A-F) after the fix.
$data = @(
    @(201,  9, 30, 15, 45, 30),
    @(1203, 3, 15, 15, 15, 15),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6, 45, 30, 60, 45),
    @(501,  9, 52, 30, 75, 45),
    @(901, 16, 15, 45, 60, 60),
    @(401,  9, 48, 67, 75, 45),
    @(701,  3, 90, 45, 97, 15),
    @(801,  3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(101,  9, 30, 15, 60, 15),
    @(1001,18, 30, 75, 60, 72),
    @(601,  9, 60, 67, 60, 42),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(1,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0, 15, 30, 30,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowIndex = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}
